$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing A:D data to B:E.
$ws.Columns("A:A").Insert()

# New column A width (old A was 6.596372, reused by shifted B column; new A is wider for labels).
$ws.Columns("A:A").ColumnWidth = 53.5

# Row labels that belong with the shifted data rows (now rows 2-19).
$labels = @(
    "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
    "CyclomaticComplexity(CC) & EffortToImplement",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbOperands & NbOperands",
    "NbOperands & EffortToImplement",
    "NbUniqueOperators & NbUniqueOperators",
    "NbOperators & NbOperators",
    "NbOperators & EffortToImplement",
    "ProgramLength & ProgramLength",
    "VocabularySize & VocabularySize",
    "ProgramVolume & ProgramVolume",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & CyclomaticComplexity(CC)",
    "EffortToImplement & NbOperands",
    "EffortToImplement & NbOperators",
    "EffortToImplement & EffortToImplement",
    "TimeToImplement & TimeToImplement"
)

# Push the existing data rows down one row (old row 1 -> new row 2, ... old row 18 -> new row 19)
# to make room for the new header row.
$ws.Rows("1:1").Insert()

# Header row for the now-shifted numeric columns (B:E).
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
